# Weekly update: a new price-report row is inserted at row 228 (pushing the
# existing rows 228-262 down to 229-263), and the new row 228 is populated
# with this week's observation for Feria Lagunitas de Puerto Montt - Plátano.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 228; formatting (incl. the date
# number format on column D) is carried along automatically.
$ws.Rows.Item(228).Insert()

$ws.Cells.Item(228, 1).Value = 4
$ws.Cells.Item(228, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(228, 3).Value = "Los Lagos"
$ws.Cells.Item(228, 4).Value = 44491
$ws.Cells.Item(228, 5).Value = 10
$ws.Cells.Item(228, 6).Value = "Fruta"
$ws.Cells.Item(228, 7).Value = 100108
$ws.Cells.Item(228, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(228, 9).Value = 100108006
$ws.Cells.Item(228, 10).Value = "Plátano"
$ws.Cells.Item(228, 11).Value = "Sin especificar"
$ws.Cells.Item(228, 12).Value = "Primera Pintón"
$ws.Cells.Item(228, 13).Value = 1400
$ws.Cells.Item(228, 14).Value = 27000
$ws.Cells.Item(228, 15).Value = 28000
$ws.Cells.Item(228, 16).Value = 27500
$ws.Cells.Item(228, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(228, 18).Value = "Ecuador"
$ws.Cells.Item(228, 19).Value = 1375
$ws.Cells.Item(228, 20).Value = 20
